$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "Electric,HVAC,Plumbing"
$ws.Range("J3").Value = "Electric,HVAC"
$ws.Range("J4").Value = "Plumbing,heating"

$ws.Range("Q5").Value = "Electric,HVAC,Plumbing"
$ws.Range("Q6").Value = "Electric,HVAC"
$ws.Range("Q7").Value = "Plumbing,heating"
